$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2318059299191375
$ws.Range("C2").Value = 0.4905660377358491
$ws.Range("J2").Value = 0.01617250673854448
$ws.Range("P2").Value = 0.1617250673854447
$ws.Range("S2").Value = 0.09973045822102426
$ws.Range("B3").Value = 0.005291005291005291
$ws.Range("C3").Value = 0.02645502645502645
$ws.Range("J3").Value = 0.01587301587301587
$ws.Range("P3").Value = 0.6613756613756614
$ws.Range("S3").Value = 0.291005291005291
$ws.Range("J4").Value = 0.07142857142857142
$ws.Range("P4").Value = 0.4464285714285715
$ws.Range("S4").Value = 0.4821428571428572
$ws.Range("B6").Value = 0.05909090909090909
$ws.Range("D6").Value = 0.00909090909090909
$ws.Range("F6").Value = 0.05
$ws.Range("J6").Value = 0.2454545454545455
$ws.Range("O6").Value = 0.01818181818181818
$ws.Range("Q6").Value = 0.2
$ws.Range("R6").Value = 0.04545454545454546
$ws.Range("S6").Value = 0.3727272727272727
$ws.Range("B7").Value = 0.1206030150753769
$ws.Range("D7").Value = 0.03015075376884422
$ws.Range("F7").Value = 0.04020100502512563
$ws.Range("J7").Value = 0.1256281407035176
$ws.Range("O7").Value = 0.01507537688442211
$ws.Range("Q7").Value = 0.1708542713567839
$ws.Range("R7").Value = 0.07537688442211055
$ws.Range("S7").Value = 0.4221105527638191
$ws.Range("B8").Value = 0.1026856240126382
$ws.Range("D8").Value = 0.01895734597156398
$ws.Range("F8").Value = 0.06477093206951026
$ws.Range("J8").Value = 0.1406003159557662
$ws.Range("O8").Value = 0.01421800947867299
$ws.Range("Q8").Value = 0.1627172195892575
$ws.Range("R8").Value = 0.0679304897314376
$ws.Range("S8").Value = 0.4281200631911533
$ws.Range("B9").Value = 0.09523809523809523
$ws.Range("D9").Value = 0.01904761904761905
$ws.Range("F9").Value = 0.05714285714285714
$ws.Range("J9").Value = 0.1523809523809524
$ws.Range("O9").Value = 0.01428571428571429
$ws.Range("Q9").Value = 0.1333333333333333
$ws.Range("R9").Value = 0.07142857142857142
$ws.Range("S9").Value = 0.4571428571428571
$ws.Range("B10").Value = 0.1165553080920564
$ws.Range("D10").Value = 0.0244988864142539
$ws.Range("E10").Value = 0.0007423904974016332
$ws.Range("F10").Value = 0.06087602078693393
$ws.Range("J10").Value = 0.1202672605790646
$ws.Range("O10").Value = 0.009651076466221232
$ws.Range("Q10").Value = 0.2487008166295472
$ws.Range("R10").Value = 0.04305864884929473
$ws.Range("S10").Value = 0.3756495916852264
$ws.Range("G11").Value = 0.1533546325878594
$ws.Range("J11").Value = 0.08945686900958466
$ws.Range("K11").Value = 0.2108626198083067
$ws.Range("L11").Value = 0.5303514376996805
$ws.Range("S11").Value = 0.01597444089456869
$ws.Range("G12").Value = 0.7485029940119761
$ws.Range("J12").Value = 0.1976047904191617
$ws.Range("K12").Value = 0.005988023952095809
$ws.Range("L12").Value = 0.02395209580838323
$ws.Range("S12").Value = 0.02395209580838323
$ws.Range("G13").Value = 0.8
$ws.Range("J13").Value = 0.15
$ws.Range("S13").Value = 0.05
$ws.Range("F15").Value = 0.02469135802469136
$ws.Range("H15").Value = 0.1790123456790123
$ws.Range("I15").Value = 0.08641975308641975
$ws.Range("J15").Value = 0.3271604938271605
$ws.Range("K15").Value = 0.06790123456790123
$ws.Range("M15").Value = 0.01851851851851852
$ws.Range("N15").Value = 0.006172839506172839
$ws.Range("O15").Value = 0.03703703703703703
$ws.Range("S15").Value = 0.2530864197530864
$ws.Range("F16").Value = 0.0196078431372549
$ws.Range("H16").Value = 0.1764705882352941
$ws.Range("I16").Value = 0.08823529411764706
$ws.Range("J16").Value = 0.4068627450980392
$ws.Range("K16").Value = 0.06862745098039216
$ws.Range("M16").Value = 0.0196078431372549
$ws.Range("O16").Value = 0.03431372549019608
$ws.Range("S16").Value = 0.1862745098039216
$ws.Range("F17").Value = 0.01289134438305709
$ws.Range("H17").Value = 0.2467771639042357
$ws.Range("I17").Value = 0.09023941068139964
$ws.Range("J17").Value = 0.4088397790055249
$ws.Range("K17").Value = 0.07734806629834254
$ws.Range("M17").Value = 0.01289134438305709
$ws.Range("O17").Value = 0.03314917127071823
$ws.Range("S17").Value = 0.1178637200736648
$ws.Range("F18").Value = 0.01418439716312057
$ws.Range("H18").Value = 0.1702127659574468
$ws.Range("I18").Value = 0.09929078014184398
$ws.Range("J18").Value = 0.4822695035460993
$ws.Range("K18").Value = 0.06382978723404255
$ws.Range("M18").Value = 0.02127659574468085
$ws.Range("O18").Value = 0.02836879432624113
$ws.Range("S18").Value = 0.1205673758865248
$ws.Range("F19").Value = 0.01052631578947368
$ws.Range("H19").Value = 0.2717105263157895
$ws.Range("I19").Value = 0.07763157894736843
$ws.Range("J19").Value = 0.3282894736842105
$ws.Range("K19").Value = 0.1105263157894737
$ws.Range("M19").Value = 0.01842105263157895
$ws.Range("N19").Value = 0.0006578947368421052
$ws.Range("O19").Value = 0.04736842105263158
$ws.Range("S19").Value = 0.1348684210526316
